$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3295464829134293
$ws.Range("C2").Value = 0.3315453850542985
$ws.Range("B3").Value = 37.7361483407898
$ws.Range("C3").Value = 37.50639873651109
$ws.Range("B4").Value = 726.4443146190338
$ws.Range("C4").Value = 704.6379668677685
$ws.Range("B5").Value = 64.97669061914591
$ws.Range("C5").Value = 61.91421887264984
$ws.Range("B6").Value = 21448.10425557612
$ws.Range("C6").Value = 19736.43690544223
$ws.Range("B7").Value = -0.006937442698816025
$ws.Range("C7").Value = -0.01842710543245034
$ws.Range("B8").Value = 1588.221000466593
$ws.Range("C8").Value = 1538.056307888663
$ws.Range("B9").Value = 1873.314312147222
$ws.Range("C9").Value = 1831.952194209138
$ws.Range("B10").Value = -0.00762233167811645
$ws.Range("C10").Value = -0.01990686176078774
$ws.Range("B11").Value = 1743.180212615688
$ws.Range("C11").Value = 1671.817174129244
$ws.Range("B12").Value = -3.938369373648492
$ws.Range("C12").Value = -3.938136514261241
$ws.Range("B13").Value = -1.586598479324239
$ws.Range("C13").Value = -1.691806061379883
$ws.Range("B14").Value = -1.533741592974951
$ws.Range("C14").Value = -1.638036194719044
$ws.Range("B15").Value = 1.827699465915042
$ws.Range("C15").Value = 1.779706249024903
